$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Burndown sheet: Sprint 3 summary row (row 5) now pulls real totals via
# formulas instead of the old hard-coded placeholder numbers.
# ---------------------------------------------------------------------------
$burndown = $wb.Worksheets.Item("Burndown")
$burndown.Range("E5").Formula = "=SUM(Sprint3!G2:G27)"
$burndown.Range("F5").Formula = "=SUM(Sprint3!H2:H27)"
$burndown.Range("G5").Copy()
$burndown.Range("F5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Backlog: US34 ("List large age differences") now has an owner assigned.
# ---------------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("D29").Value = "mb"

# ---------------------------------------------------------------------------
# Sprint3: stories are all marked Done, the last story (US34) now has
# recorded actual size/time + completed flag, and a review-results /
# keep-doing / avoid retro section was added at the bottom.
# ---------------------------------------------------------------------------
$sprint3 = $wb.Worksheets.Item("Sprint3")
$sprint3.Range("D2:D9").Value = "Done"
$sprint3.Range("G2").Value = 7
$sprint3.Range("H2").Value = 10
$sprint3.Range("I2").Value = "yes"

$sprint2 = $wb.Worksheets.Item("Sprint2")
$sprint2.Range("B28:B36").Copy()
$sprint3.Range("B12:B20").PasteSpecial(-4122)

$sprint3.Range("B12").Value = "Review Results"
$sprint3.Range("B14").Value = "Keep doing:"
$sprint3.Range("B15").Value = "Get work done ahead of time so there is time for integration testing"
$sprint3.Range("B18").Value = "Avoid:"
$sprint3.Range("B19").Value = "Lack of communication"
$sprint3.Range("B20").Value = "Push regularly"

# ---------------------------------------------------------------------------
# Sprint4: add the US34 entry for this sprint as well.
# ---------------------------------------------------------------------------
$sprint4 = $wb.Worksheets.Item("Sprint4")
$sprint4.Range("C4").Copy()
$sprint4.Range("C8").PasteSpecial(-4122)
$sprint4.Range("A8").Value = "US34"
$sprint4.Range("B8").Value = "List large age differences"
$sprint4.Range("C8").Value = "mb"
$sprint4.Range("E8").Value = 8
$sprint4.Range("F8").Value = 15

# ---------------------------------------------------------------------------
# Selections on each sheet (restored / updated to match where the authors
# were last working), and which sheet/tab is active.
# ---------------------------------------------------------------------------
$backlog.Range("A29:D29").Select()
$burndown.Range("F6").Select()
$sprint2.Range("B28:B36").Select()
$sprint3.Range("D15").Select()
$sprint4.Range("E9").Select()

$stories = $wb.Worksheets.Item("Stories")
$stories.Range("D40").Select()

$team = $wb.Worksheets.Item("Team")
$team.Activate()
$team.Range("D24").Select()
